$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header: "sup." -> "man. pn."
$ws.Range("C1").Value = "man. pn."

# 2. Remove D6 ("Need Model") entirely
$ws.Range("D6").ClearContents()

# 3. F1 row: supplier part number changes
$ws.Range("C10").Value = "023403.5MXP"

# 4. New row 11: C11 part number changes, A11 gets "F1 Holder"
$ws.Range("C11").Value = "05200101Z"
$ws.Range("A11").Value = "F1 Holder"

# 5. TP1/TP2/TP3 row: C12 becomes numeric 5271 (was text "36-5271-ND"),
#    left-aligned like the other numeric hyperlink cell (C8), and the
#    hyperlink keeps its display text in sync.
$ws.Range("C12").Value = 5271
$ws.Range("C12").HorizontalAlignment = -4131

foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$12') {
        $hl.TextToDisplay = "5271"
    }
}

# 6. Selection moves to A12
[void]$ws.Range("A12").Select()
